# Apply the refreshed crypto-price snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.564.96'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '3.208.99'
$ws.Range("E3").Value = '  +4.51%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'" + '239.78'
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("D6").Value = "'" + '619.66'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  +5.39%  '
$ws.Range("D8").Value = "'" + '0.371'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").Value = "'" + '1.00'
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").Value = '3.194.09'
$ws.Range("E10").Value = '  +4.07%  '
$ws.Range("D11").Value = "'" + '0.735'
$ws.Range("E11").Value = '  +4.83%  '
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("E13").Value = '  -0.23%  '
$ws.Range("D14").Value = "'" + '35.36'
$ws.Range("E14").Value = '  +1.92%  '
$ws.Range("E15").Value = '  +3.73%  '
$ws.Range("D16").Value = '90.793.58'
$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("D17").Value = '3.759.86'
$ws.Range("E17").Value = '  +3.36%  '
$ws.Range("D18").Value = '3.201.09'
$ws.Range("E18").Value = '  +4.74%  '
$ws.Range("D19").Value = "'" + '3.70'
$ws.Range("E19").Value = '  -2.72%  '
$ws.Range("D20").Value = "'" + '15.20'
$ws.Range("E20").Value = '  +10.58%  '
$ws.Range("D21").Value = "'" + '6.03'
$ws.Range("E21").Value = '  +11.49%  '
$ws.Range("D22").Value = "'" + '451.70'
$ws.Range("E22").Value = '  +4.73%  '
$ws.Range("E23").Value = '  -4.23%  '
$ws.Range("D24").Value = "'" + '9.23'
$ws.Range("E24").Value = '  +5.98%  '
$ws.Range("D25").Value = "'" + '5.78'
$ws.Range("E25").Value = '  +4.07%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = "'" + '89.11'
$ws.Range("E26").Value = '  +2.54%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").Value = "'" + '12.00'
$ws.Range("E27").Value = '  +2.38%  '
$ws.Range("E28").Value = '  +3.42%  '
$ws.Range("D29").Value = "'" + '1.00'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("E30").Value = '  +56.02%  '
$ws.Range("D31").Value = "'" + '0.234'
$ws.Range("E31").Value = '  +18.38%  '
$ws.Range("E32").Value = '  +8.09%  '
$ws.Range("D33").Value = "'" + '9.39'
$ws.Range("E33").Value = '  +4.83%  '
$ws.Range("E34").Value = '  +14.49%  '
$ws.Range("E35").Value = '  -5.97%  '
$ws.Range("D36").Value = "'" + '27.30'
$ws.Range("E36").Value = '  +6.60%  '
$ws.Range("D37").Value = "'" + '7.71'
$ws.Range("E37").Value = '  +8.00%  '
$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").Value = "'" + '1.98'
$ws.Range("E38").Value = '  +5.83%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = "'" + '510.16'
$ws.Range("E39").Value = '  +3.91%  '
$ws.Range("E40").Value = '  +7.67%  '
$ws.Range("D41").Value = "'" + '0.455'
$ws.Range("E41").Value = '  +13.96%  '
$ws.Range("D42").Value = "'" + '3.81'
$ws.Range("E42").Value = '  -10.21%  '
$ws.Range("D43").Value = "'" + '3.44'
$ws.Range("E43").Value = '  -4.63%  '
$ws.Range("D44").Value = "'" + '22.07'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D46").Value = "'" + '0.739'
$ws.Range("E46").Value = '  +9.93%  '
$ws.Range("D47").Value = "'" + '1.93'
$ws.Range("E47").Value = '  +4.35%  '
$ws.Range("D48").Value = "'" + '156.38'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").Value = "'" + '1.39'
$ws.Range("E49").Value = '  +7.23%  '
$ws.Range("D50").Value = "'" + '4.48'
$ws.Range("E50").Value = '  +3.80%  '
$ws.Range("D51").Value = "'" + '43.95'
$ws.Range("E51").Value = '  -0.85%  '
